# Apply the authored edits to the "Dictionnaire de données" workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from the default "Feuil1" to the descriptive French title.
$ws.Name = "Dictionnaire de données"

# Tighten up the description text for the "NumeroPermis" field (row 11, col B).
$ws.Range("B11").Value = "Permis de l'inventeur"

# Zoom out from 145% to 115% on the sheet's window.
$excel.ActiveWindow.Zoom = 115

# Select the whole of row 11 (as if the row header had been clicked),
# moving the active cell from D11 to A11.
$ws.Rows(11).Select()
